# "switch from hid to cdc-acm"
#
# The "message IDs" table on Sheet1 (columns I:M, rows 34-65) loses the
# HID-specific "0x03 / define USB frame size / U8 / definitely send this
# one before subscribing to anything" entry that used to live on row 34.
# Every following entry (rows 36-65) shifts up by one row (into rows
# 35-64), and the now-unused last row (65) disappears.
#
# We reproduce that by copying each row's I:M content up from the row
# below, cell-by-cell (starting at the top so we never overwrite data we
# still need), then clearing the now-surplus last row.  Column C/D (and
# A/B where present) are untouched - only I:M move.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$firstRow = 34
$lastRow = 65
$cols = @("I", "J", "K", "L", "M")

foreach ($col in $cols) {
    for ($r = $firstRow; $r -lt $lastRow; $r++) {
        $srcCell = $ws.Range($col + ($r + 1))
        $dstCell = $ws.Range($col + $r)

        # Clear the destination first: Copy-ing an empty source cell does
        # not blank out a previously-populated destination on its own.
        $dstCell.ClearContents()

        if ($srcCell.Value() -ne $null) {
            $srcCell.Copy($dstCell)
        }
    }

    # The last row had nothing shifted into it - make sure it's empty.
    $ws.Range($col + $lastRow).ClearContents()
}

# Restore the view state recorded in the saved file: scrolled down/right a
# bit further, with N34 as the active cell.
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("N34").Select()
